$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.785.82"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = "'2.115.10"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.27%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'332.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.35%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = "'0.5329"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.54%  '
$ws.Range('D8').Value = "'0.4398"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.99%  '
$ws.Range('D9').Value = "'0.09003"
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Value = "'47.37"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.96%  '
$ws.Range('D11').Value = "'1.180"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.38%  '
$ws.Range('D12').Value = "'24.96"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.41%  '
$ws.Range('D13').Value = "'2.109.09"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.85%  '
$ws.Range('D14').Value = "'6.766"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.38%  '
$ws.Range('D15').Value = "'7.801"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.52%  '
$ws.Range('D16').Value = "'96.68"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Value = "'0.00001131"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').Value = "'0.06679"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').Value = "'19.14"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').Value = "'1.000"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'6.334"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.31%  '
$ws.Range('D23').Value = "'30.853.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.84%  '
$ws.Range('D24').Value = "'12.34"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.43%  '
$ws.Range('D25').Value = "'2.358.25"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.17%  '
$ws.Range('D26').Value = "'2.284"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('D27').Value = "'22.80"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = "'2.587"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.81%  '
$ws.Range('D29').Value = "'163.01"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = "'133.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('E31').Value = '  +4.64%  '
$ws.Range('D32').Value = "'0.1083"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('D33').Value = "'6.241"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.08%  '
$ws.Range('D34').Value = "'4.010"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.22%  '
$ws.Range('D35').Value = "'1.562"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +18.90%  '
$ws.Range('D36').Value = "'0.02599"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.82%  '
$ws.Range('D37').Value = "'12.90"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.61%  '
$ws.Range('D38').Value = "'5.531"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('D39').Value = "'0.06755"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.69%  '
$ws.Range('D40').Value = "'9.490"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.71%  '
$ws.Range('E41').Value = '  +4.97%  '
$ws.Range('D42').Value = "'0.6837"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.19%  '
$ws.Range('D43').Value = "'1.250"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'14.16"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.6448"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.32%  '
$ws.Range('D46').Value = "'0.9994"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').Value = "'2.231"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('D48').Value = "'3.659"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').Value = "'1.265"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.35%  '
$ws.Range('D50').Value = "'82.98"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.45%  '
$ws.Range('D51').Value = "'121.59"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.89%  '
